# Add a new "Ev_kW" column (electric-vehicle load, kW) to the
# INTERNAL_LOADS sheet, right after the existing "Qcpro_Wm2" column (M).
# Header goes in N1, and every data row (2-20) gets a default value of 0,
# matching the formatting already used for the other numeric load columns.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("INDOOR_COMFORT")
$ws2 = $wb.Worksheets.Item("INTERNAL_LOADS")

# Copy the formatting of column M (the last existing column) into the new
# column N so the header/cell styling (fill, border, number format, etc.)
# stays consistent with the rest of the table.
$ws2.Range("M1").Copy() | Out-Null
$ws2.Range("N1").PasteSpecial(-4122) | Out-Null
$ws2.Range("N1").Value = "Ev_kW"

for ($r = 2; $r -le 20; $r++) {
    $ws2.Range("M$r").Copy() | Out-Null
    $ws2.Range("N$r").PasteSpecial(-4122) | Out-Null
    $ws2.Cells.Item($r, 14).Value = 0
}

$excel.CutCopyMode = $false

# Restore the last-active-cell selections seen in the authored workbook.
$ws1.Range("D34").Select() | Out-Null
$ws2.Range("N1").Select() | Out-Null
